$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells being updated so that
# numeric-looking strings (e.g. "2.00", "0.0850") keep their exact text
# formatting instead of being coerced into numbers.
$priceCells = @("D2","D3","D5","D6","D7","D9","D10","D11","D15","D16","D18","D19","D21","D22","D23","D24","D26","D30","D31","D32","D33","D34","D35","D38","D39","D40","D42","D43","D44","D45","D46","D48","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range('D2').Value = '51.829.25'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '2.779.69'
$ws.Range('E3').Value = '  -1.82%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '357.48'
$ws.Range('E5').Value = '  +1.31%  '
$ws.Range('D6').Value = '109.75'
$ws.Range('E6').Value = '  -3.00%  '
$ws.Range('D7').Value = '0.566'
$ws.Range('E7').Value = '  +0.97%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '0.595'
$ws.Range('E9').Value = '  -0.81%  '
$ws.Range('D10').Value = '40.04'
$ws.Range('E10').Value = '  -3.58%  '
$ws.Range('D11').Value = '0.0850'
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('E12').Value = '  +0.74%  '
$ws.Range('E13').Value = '  -2.38%  '
$ws.Range('E14').Value = '  -1.44%  '
$ws.Range('D15').Value = '3.215.18'
$ws.Range('E15').Value = '  -1.93%  '
$ws.Range('D16').Value = '2.781.41'
$ws.Range('E16').Value = '  -1.81%  '
$ws.Range('E17').Value = '  +4.45%  '
$ws.Range('D18').Value = '51.768.97'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('D19').Value = '7.42'
$ws.Range('E19').Value = '  +0.57%  '
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('D21').Value = '13.06'
$ws.Range('E21').Value = '  -2.86%  '
$ws.Range('D22').Value = '0.0₃0977'
$ws.Range('E22').Value = '  -1.58%  '
$ws.Range('D23').Value = '273.47'
$ws.Range('E23').Value = '  +1.10%  '
$ws.Range('D24').Value = '70.07'
$ws.Range('E24').Value = '  +0.60%  '
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('D26').Value = '26.64'
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('E28').Value = '  -1.09%  '
$ws.Range('E29').Value = '  +3.96%  '
$ws.Range('D30').Value = '2.21'
$ws.Range('E30').Value = '  -1.58%  '
$ws.Range('D31').Value = '0.0464'
$ws.Range('E31').Value = '  +4.00%  '
$ws.Range('D32').Value = '51.52'
$ws.Range('E32').Value = '  +1.68%  '
$ws.Range('D33').Value = '33.87'
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').Value = '5.69'
$ws.Range('E34').Value = '  -1.87%  '
$ws.Range('D35').Value = '0.0844'
$ws.Range('E35').Value = '  +2.32%  '
$ws.Range('E36').Value = '  +7.88%  '
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('D38').Value = '3.24'
$ws.Range('E38').Value = '  +1.02%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '2.00'
$ws.Range('E39').Value = '  -3.61%  '
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').Value = '18.07'
$ws.Range('E40').Value = '  +0.29%  '
$ws.Range('E41').Value = '  -0.36%  '
$ws.Range('D42').Value = '2.52'
$ws.Range('E42').Value = '  -1.33%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = '2.24'
$ws.Range('E43').Value = '  -2.59%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').Value = '121.62'
$ws.Range('E44').Value = '  -3.02%  '
$ws.Range('D45').Value = '21.97'
$ws.Range('E45').Value = '  -6.51%  '
$ws.Range('D46').Value = '2.066.24'
$ws.Range('E46').Value = '  -0.62%  '
$ws.Range('E47').Value = '  -2.09%  '
$ws.Range('D48').Value = '2.17'
$ws.Range('E48').Value = '  -6.01%  '
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('D50').Value = '0.933'
$ws.Range('E50').Value = '  -0.21%  '
$ws.Range('D51').Value = '8.93'
$ws.Range('E51').Value = '  +0.19%  '
